$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 98 (shifts old rows 98-114 down to 102-118)
$ws.Rows("98:101").Insert()

# Fill in the new rows (98-101) with this week's data for
# Agricola del Norte S.A. de Arica - Frutilla, date 2023-09-15 (serial 45173)
$newRows = @(
    @{ Row = 98;  L = "Especial"; M = 65;  N = 8000; O = 9000; P = 8462; S = 2821 },
    @{ Row = 99;  L = "Primera";  M = 100; N = 6000; O = 7000; P = 6500; S = 2167 },
    @{ Row = 100; L = "Segunda";  M = 100; N = 4000; O = 5000; P = 4500; S = 1500 },
    @{ Row = 101; L = "Tercera";  M = 100; N = 2000; O = 3000; P = 2500; S = 833 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($row, 4).Value = 45173
    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "$/bandeja 3 kilos"
    $ws.Cells.Item($row, 18).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 3
}
